$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 306
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F7").Value = 939
$ws1.Range("F10").Value = 546
$ws1.Range("F11").Value = 1406
$ws1.Range("F13").Value = 1316
$ws1.Range("F14").Value = 2973
$ws1.Range("F15").Value = 369
$ws1.Range("F16").Value = 1593
$ws1.Range("F18").Value = 774
$ws1.Range("F19").Value = 226
$ws1.Range("F20").Value = 1337
$ws1.Range("F23").Value = 1107
$ws1.Range("F24").Value = 390
$ws1.Range("F25").Value = 3411
$ws1.Range("F26").Value = 665
$ws1.Range("F28").Value = 1515

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 45
$ws2.Range("F12").Value = 69

# --- Sheet "全部类型" (sheet4, combined view, rows offset by +1 vs sheet1/2) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 306
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F13").Value = 45
$ws4.Range("F17").Value = 939
$ws4.Range("F20").Value = 546
$ws4.Range("F21").Value = 1406
$ws4.Range("F23").Value = 1316
$ws4.Range("F24").Value = 2973
$ws4.Range("F25").Value = 369
$ws4.Range("F26").Value = 1593
$ws4.Range("F28").Value = 774
$ws4.Range("F29").Value = 226
$ws4.Range("F30").Value = 1337
$ws4.Range("F35").Value = 1107
$ws4.Range("F36").Value = 390
$ws4.Range("F37").Value = 3411
$ws4.Range("F38").Value = 665
$ws4.Range("F40").Value = 1515
$ws4.Range("F41").Value = 69
